# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) previously held a different strikeout-related
# statistic ("Strike#"); it has been regenerated to hold the true strikeout
# count (K) for each outing. Update the affected cells with their new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 4
    6  = 3
    7  = 2
    8  = 2
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 0
    14 = 4
    15 = 1
    16 = 0
    17 = 0
    18 = 2
    19 = 0
    20 = 2
    21 = 0
    22 = 1
    23 = 1
    24 = 1
    25 = 0
    26 = 2
    27 = 1
    28 = 1
    29 = 1
    30 = 2
    31 = 4
    32 = 3
    33 = 2
    35 = 4
    36 = 2
    38 = 1
    39 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
